$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "29.233.55"
$ws.Range("E2").Value = "  -0.50%  "

$ws.Range("D3").Value = "1.863.23"
$ws.Range("E3").Value = "  -0.56%  "

$ws.Range("D4").Value = "'1.000"
$ws.Range("E4").Value = "  -0.01%  "

$ws.Range("D5").Value = "'0.7154"
$ws.Range("E5").Value = "  +0.20%  "

$ws.Range("D6").Value = "'240.49"
$ws.Range("E6").Value = "  +0.37%  "

$ws.Range("E7").Value = "  +0.04%  "

$ws.Range("D8").Value = "'0.3092"
$ws.Range("E8").Value = "  +0.36%  "

$ws.Range("D9").Value = "'0.07704"
$ws.Range("E9").Value = "  -1.48%  "

$ws.Range("D10").Value = "'25.07"
$ws.Range("E10").Value = "  +1.17%  "

$ws.Range("D11").Value = "'0.08328"
$ws.Range("E11").Value = "  +0.94%  "

$ws.Range("D12").Value = "1.973.32"
$ws.Range("E12").Value = "  +5.22%  "

$ws.Range("D13").Value = "'0.7176"
$ws.Range("E13").Value = "  -0.86%  "

$ws.Range("D14").Value = "'5.217"
$ws.Range("E14").Value = "  -0.82%  "

$ws.Range("D15").Value = "'90.82"
$ws.Range("E15").Value = "  -0.47%  "

$ws.Range("D16").Value = "29.361.05"
$ws.Range("E16").Value = "  -0.08%  "

$ws.Range("B18").Value = "WrappedliquidstakedEther2.0"
$ws.Range("C18").Value = "https://coinranking.com/coin/CiixT63n3+wrappedliquidstakedether20-wsteth"
$ws.Range("D18").Value = "2.171.05"
$ws.Range("E18").Value = "  +2.66%  "

$ws.Range("B19").Value = "BitcoinCash"
$ws.Range("C19").Value = "https://coinranking.com/coin/ZlZpzOJo43mIo+bitcoincash-bch"
$ws.Range("D19").Value = "'243.67"
$ws.Range("E19").Value = "  -0.15%  "

$ws.Range("D20").Value = "'0.000007802"
$ws.Range("E20").Value = "  -1.25%  "

$ws.Range("E21").Value = "  -0.83%  "

$ws.Range("B22").Value = "Dai"
$ws.Range("C22").Value = "https://coinranking.com/coin/MoTuySvg7+dai-dai"
$ws.Range("D22").Value = "'1.001"
$ws.Range("E22").Value = "  +0.16%  "

$ws.Range("B23").Value = "Chainlink"
$ws.Range("C23").Value = "https://coinranking.com/coin/VLqpJwogdhHNb+chainlink-link"
$ws.Range("D23").Value = "'8.003"
$ws.Range("E23").Value = "  +0.75%  "

$ws.Range("E24").Value = "  +0.05%  "

$ws.Range("D25").Value = "'0.1612"
$ws.Range("E25").Value = "  +3.88%  "

$ws.Range("D26").Value = "'162.89"
$ws.Range("E26").Value = "  -0.39%  "

$ws.Range("D27").Value = "'8.911"
$ws.Range("E27").Value = "  -0.85%  "

$ws.Range("D28").Value = "'18.59"
$ws.Range("E28").Value = "  +1.55%  "

$ws.Range("D29").Value = "'1.342"
$ws.Range("E29").Value = "  -1.19%  "

$ws.Range("D30").Value = "'4.443"
$ws.Range("E30").Value = "  +1.69%  "

$ws.Range("D31").Value = "'1.496"
$ws.Range("E31").Value = "  +0.93%  "

$ws.Range("E32").Value = "  +3.49%  "

$ws.Range("D33").Value = "'0.05188"
$ws.Range("E33").Value = "  -1.47%  "

$ws.Range("B34").Value = "ImmutableX"
$ws.Range("C34").Value = "https://coinranking.com/coin/Z96jIvLU7+immutablex-imx"
$ws.Range("D34").Value = "'0.7957"
$ws.Range("E34").Value = "  +10.83%  "

$ws.Range("B35").Value = "LidoDAOToken"
$ws.Range("C35").Value = "https://coinranking.com/coin/Pe93bIOD2+lidodaotoken-ldo"
$ws.Range("D35").Value = "'1.926"
$ws.Range("E35").Value = "  +0.35%  "

$ws.Range("B36").Value = "ARBITRUM"
$ws.Range("C36").Value = "https://coinranking.com/coin/1Uo6s62Oc+arbitrum-arb"
$ws.Range("D36").Value = "'1.173"
$ws.Range("E36").Value = "  -2.00%  "

$ws.Range("D37").Value = "'2.683"
$ws.Range("E37").Value = "  +0.10%  "

$ws.Range("D38").Value = "'0.01857"
$ws.Range("E38").Value = "  -0.04%  "

$ws.Range("D39").Value = "'2.698"
$ws.Range("E39").Value = "  -0.37%  "

$ws.Range("D40").Value = "1.183.29"
$ws.Range("E40").Value = "  -1.84%  "

$ws.Range("D41").Value = "'6.264"
$ws.Range("E41").Value = "  +3.40%  "

$ws.Range("D42").Value = "'0.9042"
$ws.Range("E42").Value = "  -0.25%  "

$ws.Range("D43").Value = "'73.11"
$ws.Range("E43").Value = "  +1.12%  "

$ws.Range("E44").Value = "  -0.03%  "

$ws.Range("B45").Value = "RocketPoolETH"
$ws.Range("C45").Value = "https://coinranking.com/coin/QJZRUGyNI+rocketpooleth-reth"
$ws.Range("D45").Value = "2.066.74"
$ws.Range("E45").Value = "  +2.54%  "

$ws.Range("B46").Value = "Quant"
$ws.Range("C46").Value = "https://coinranking.com/coin/bauj_21eYVwso+quant-qnt"
$ws.Range("D46").Value = "'102.47"
$ws.Range("E46").Value = "  -0.74%  "

$ws.Range("D47").Value = "'0.5207"
$ws.Range("E47").Value = "  -2.50%  "

$ws.Range("D49").Value = "'9.352"
$ws.Range("E49").Value = "  +1.37%  "

$ws.Range("D50").Value = "'1.014"
$ws.Range("E50").Value = "  +1.45%  "

$ws.Range("B51").Value = "Aptos"
$ws.Range("C51").Value = "https://coinranking.com/coin/HGYj5JCv5+aptos-apt"
$ws.Range("D51").Value = "'7.080"
$ws.Range("E51").Value = "  +1.09%  "
